# Add a new "2022-Q3" sheet (duplicate of "2022-Q2" layout) in front of
# "2022-Q2", fill it with the new quarter's fund data, and insert a
# matching summary row at the top of the "总计" sheet's data block.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by duplicating "2022-Q2" (keeps the
#    same column layout/styles), placed immediately before "2022-Q2".
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Overwrite the copied rows with the 2022-Q3 fund holdings.
$q3.Cells.Item(2,2).Value = "'501029"
$q3.Cells.Item(2,3).Value = "华宝标普中国A股红利机会指数（LOF）A"
$q3.Cells.Item(2,4).Value = "'8.11"
$q3.Cells.Item(2,5).Value = "'94.26"
$q3.Cells.Item(2,6).Value = "'1.62"
$q3.Cells.Item(2,7).Value = "'0.1314"
$q3.Cells.Item(2,8).Value = 5

$q3.Cells.Item(3,2).Value = "'005125"
$q3.Cells.Item(3,3).Value = "华宝标普中国A股红利机会指数C"
$q3.Cells.Item(3,4).Value = "'3.38"
$q3.Cells.Item(3,5).Value = "'94.26"
$q3.Cells.Item(3,6).Value = "'1.62"
$q3.Cells.Item(3,7).Value = "'0.0548"
$q3.Cells.Item(3,8).Value = 5

# ---------------------------------------------------------------------
# 2) Insert a new row 2 into "总计" for the 2022-Q3 summary, shifting the
#    existing quarters down by one row, and renumber the index column.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# The freshly-inserted row doesn't inherit the surrounding cell styles by
# itself, so copy them over from the row just below (which still holds
# the original "2022-Q2" formatting) before filling in the new values.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B3:D3").Copy()
$total.Range("B2:D2").PasteSpecial(-4122)

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 2
$total.Cells.Item(2,4).Value = 0.19

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4

# Restore "2021-Q3" (last tab) as the active/selected sheet — copying and
# renaming the new "2022-Q3" sheet above shifted the active tab onto it.
$wb.Worksheets.Item("2021-Q3").Activate()
